$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.787096774193548
$ws.Range("C2").Value = 0.71334214002642
$ws.Range("D2").Value = 0.821362799263352
$ws.Range("E2").Value = 0.725155279503106
$ws.Range("F2").Value = 0.53424139235717

$ws.Range("B3").Value = 0.75
$ws.Range("C3").Value = 0.667107001321004
$ws.Range("D3").Value = 0.801104972375691
$ws.Range("E3").Value = 0.740683229813665
$ws.Range("F3").Value = 0.552402572833901

$ws.Range("B4").Value = 0.783870967741935
$ws.Range("C4").Value = 0.690885072655218
$ws.Range("D4").Value = 0.815837937384899
$ws.Range("E4").Value = 0.680124223602484
$ws.Range("F4").Value = 0.52099886492622
